$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5; existing rows 5.. shift down to 6..
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with its data
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = "Vega Monumental Concepción"
$ws.Range("C5").Value = "Bíobío"
$ws.Range("D5").Value = 44812
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100108
$ws.Range("H5").Value = "Tropicales y subtropicales"
$ws.Range("I5").Value = 100108002
$ws.Range("J5").Value = "Mango"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 140
$ws.Range("N5").Value = 8500
$ws.Range("O5").Value = 9000
$ws.Range("P5").Value = 8714
$ws.Range("Q5").Value = "$/bandeja 4 kilos"
$ws.Range("R5").Value = "Perú"
$ws.Range("S5").Value = 2178
$ws.Range("T5").Value = 4
